$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.741.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.885.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.52%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4731'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3988'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.92'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08070'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.78%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.85'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.899.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.979'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.214'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.05'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001042'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06602'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.27'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.744.93'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.526'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.106.12'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.22'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.81%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.102'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.594'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.59'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9706'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09546'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.471'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.625'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.306'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06131'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02259'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.231'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.177'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6011'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1903'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.262'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5703'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.26'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.409'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06823'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.69'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.36%  '
